$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 692469.2
$ws.Range("I80").Value = 933
$ws.Range("J80").Value = 865353.25
$ws.Range("K80").Value = 2799
$ws.Range("L80").Value = 2596059.75
$ws.Range("M80").Value = -1801
$ws.Range("N80").Value = -2598055.75
$ws.Range("H83").Value = 692469.2
$ws.Range("I83").Value = 933
$ws.Range("J83").Value = 865353.25
$ws.Range("K83").Value = 8397
$ws.Range("L83").Value = 7788179.25
$ws.Range("M83").Value = -3405
$ws.Range("N83").Value = -7798163.25
$ws.Range("H92").Value = 422.7
$ws.Range("I92").Value = 392.3158
$ws.Range("J92").Value = 1000
$ws.Range("K92").Value = 392.3158
$ws.Range("L92").Value = 1000
$ws.Range("M92").Value = 855.6841999999999
$ws.Range("N92").Value = -3496
$ws.Range("H138").Value = 4763.6895
$ws.Range("I138").Value = 2008.2258
$ws.Range("J138").Value = 7927.3706
$ws.Range("K138").Value = 6024.6774
$ws.Range("L138").Value = 23782.1118
$ws.Range("M138").Value = -884.6773999999996
$ws.Range("N138").Value = -34062.1118
$ws.Range("H139").Value = 29714.285
$ws.Range("J139").Value = 29714.285
$ws.Range("L139").Value = 29714.285
$ws.Range("N139").Value = -39994.285
$ws.Range("H140").Value = 24666.666
$ws.Range("J140").Value = 24666.666
$ws.Range("L140").Value = 24666.666
$ws.Range("N140").Value = -35026.666
$ws.Range("H141").Value = 719130.1
$ws.Range("I141").Value = 1611.75
$ws.Range("J141").Value = 1539151.1
$ws.Range("K141").Value = 4835.25
$ws.Range("L141").Value = 4617453.300000001
$ws.Range("M141").Value = 344.75
$ws.Range("N141").Value = -4627813.300000001
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2365.34
$ws.Range("I32").Value = 2365.34
$ws.Range("K32").Value = 2365.34
$ws.Range("M32").Value = -2078.34
$ws.Range("H74").Value = 1526.9688
$ws.Range("I74").Value = 1091.6957
$ws.Range("J74").Value = 2639.3333
$ws.Range("K74").Value = 1091.6957
$ws.Range("L74").Value = 2639.3333
$ws.Range("M74").Value = -217.6957
$ws.Range("N74").Value = -4387.3333
$ws.Range("H77").Value = 1526.9688
$ws.Range("I77").Value = 1091.6957
$ws.Range("J77").Value = 2639.3333
$ws.Range("K77").Value = 5458.4785
$ws.Range("L77").Value = 13196.6665
$ws.Range("M77").Value = -1090.4785
$ws.Range("N77").Value = -21932.6665
$ws.Range("H94").Value = 30011.818
$ws.Range("J94").Value = 30011.818
$ws.Range("L94").Value = 30011.818
$ws.Range("N94").Value = -31813.818
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2878.459
$ws.Range("I134").Value = 2697.4614
$ws.Range("K134").Value = 8092.3842
$ws.Range("M134").Value = -5557.3842
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H88").Value = 30000
$ws.Range("J88").Value = 30000
$ws.Range("L88").Value = 30000
$ws.Range("N88").Value = -30812
$ws.Range("H91").Value = 30000
$ws.Range("J91").Value = 30000
$ws.Range("L91").Value = 30000
$ws.Range("N91").Value = -32808
$ws.Range("H105").Value = 2219.0715
$ws.Range("I105").Value = 1741
$ws.Range("K105").Value = 1741
$ws.Range("M105").Value = 6
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 907.01
$ws.Range("I131").Value = 449.8
$ws.Range("J131").Value = 931.07367
$ws.Range("K131").Value = 1349.4
$ws.Range("L131").Value = 2793.22101
$ws.Range("M131").Value = 3690.6
$ws.Range("N131").Value = -12873.22101
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 3647.5
$ws.Range("I97").Value = 1530
$ws.Range("J97").Value = 10000
$ws.Range("K97").Value = 1530
$ws.Range("L97").Value = 10000
$ws.Range("M97").Value = -1034
$ws.Range("N97").Value = -10992
$ws.Range("H113").Value = 3628.1304
$ws.Range("I113").Value = 2920.4119
$ws.Range("J113").Value = 5633.3335
$ws.Range("K113").Value = 2920.4119
$ws.Range("L113").Value = 5633.3335
$ws.Range("M113").Value = -750.4119000000001
$ws.Range("N113").Value = -9973.333500000001
$ws.Range("H122").Value = 5771.724
$ws.Range("I122").Value = 5150
$ws.Range("K122").Value = 15450
$ws.Range("M122").Value = -13000
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2478.3333
$ws.Range("I100").Value = 1750
$ws.Range("J100").Value = 2842.5
$ws.Range("K100").Value = 1750
$ws.Range("L100").Value = 2842.5
$ws.Range("M100").Value = -1209
$ws.Range("N100").Value = -3924.5
$ws.Range("H122").Value = 2868.5117
$ws.Range("I122").Value = 2437.5
$ws.Range("J122").Value = 3673.0667
$ws.Range("K122").Value = 7312.5
$ws.Range("L122").Value = 11019.2001
$ws.Range("M122").Value = -4862.5
$ws.Range("N122").Value = -15919.2001
$ws.Range("H132").Value = 2827.647
$ws.Range("I132").Value = 2031.5294
$ws.Range("J132").Value = 3623.7646
$ws.Range("K132").Value = 6094.5882
$ws.Range("L132").Value = 10871.2938
$ws.Range("M132").Value = -3564.5882
$ws.Range("N132").Value = -15931.2938
$ws.Range("H136").Value = 2383671.2
$ws.Range("I136").Value = 3335679.8
$ws.Range("J136").Value = 3650.4167
$ws.Range("K136").Value = 10007039.4
$ws.Range("L136").Value = 10951.2501
$ws.Range("M136").Value = -10004489.4
$ws.Range("N136").Value = -16051.2501
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 29333.334
$ws.Range("J94").Value = 29333.334
$ws.Range("L94").Value = 29333.334
$ws.Range("N94").Value = -31135.334
$ws.Range("H100").Value = 553.3333
$ws.Range("I100").Value = 416.66666
$ws.Range("J100").Value = 690
$ws.Range("K100").Value = 833.33332
$ws.Range("L100").Value = 1380
$ws.Range("M100").Value = -292.33332
$ws.Range("N100").Value = -2462
$ws.Range("H124").Value = 23214.5
$ws.Range("J124").Value = 23214.5
$ws.Range("L124").Value = 23214.5
$ws.Range("N124").Value = -33034.5
$ws.Range("H132").Value = 378322.34
$ws.Range("I132").Value = 838891.3
$ws.Range("J132").Value = 9867.134
$ws.Range("K132").Value = 2516673.9
$ws.Range("L132").Value = 29601.402
$ws.Range("M132").Value = -2514143.9
$ws.Range("N132").Value = -34661.402
$ws.Range("H136").Value = 1903.1489
$ws.Range("I136").Value = 1605.7931
$ws.Range("J136").Value = 2382.2222
$ws.Range("K136").Value = 4817.379300000001
$ws.Range("L136").Value = 7146.6666
$ws.Range("M136").Value = -2267.379300000001
$ws.Range("N136").Value = -12246.6666
$ws.Range("H141").Value = 28681.818
$ws.Range("J141").Value = 28681.818
$ws.Range("L141").Value = 28681.818
$ws.Range("N141").Value = -39041.818
